# Split-sheet-into-workbook export formatting, applied to the "France" example output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet to match the country it represents ---
$ws.Name = "France"

# --- Tab color (Accent 5, Lighter 40%) ---
$ws.Tab.Color = 15123357

# --- Build the header style (bold white text on dark-blue fill, thin border,
#     centered horizontally and top-aligned) on a scratch cell, then stamp it
#     onto the header row A1:O1 in one shot via PasteSpecial so every header
#     cell picks up a single combined style. ---
$scratch = $ws.Range("Z1")
$scratch.Font.Bold = $true
$scratch.Font.ColorIndex = 2
$scratch.Interior.Color = 9524736
$scratch.Borders.LineStyle = 1
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160
$scratch.Copy()
$ws.Range("A1:O1").PasteSpecial(-4122)
$scratch.Clear()

# --- Highlight the "Gross Sales" column (H) with a light gray fill ---
$ws.Range("H2:H141").Interior.Color = 14277081

# --- Column widths (character units; matches the widths produced by
#     Excel's auto-fit after the header formatting above was applied) ---
$ws.Columns.Item(1).ColumnWidth = 14.83
$ws.Columns.Item(2).ColumnWidth = 6.83
$ws.Columns.Item(3).ColumnWidth = 8.17
$ws.Columns.Item(4).ColumnWidth = 10.17
$ws.Columns.Item(5).ColumnWidth = 8.83
$ws.Columns.Item(6).ColumnWidth = 17.72
$ws.Columns.Item(7).ColumnWidth = 8.5
$ws.Columns.Item(8).ColumnWidth = 9.83
$ws.Columns.Item(9).ColumnWidth = 9.17
$ws.Columns.Item(10).ColumnWidth = 10.17
$ws.Columns.Item(11).ColumnWidth = 7.17
$ws.Columns.Item(12).ColumnWidth = 9.17
$ws.Columns.Item(13).ColumnWidth = 13.61
$ws.Columns.Item(14).ColumnWidth = 11.5
$ws.Columns.Item(15).ColumnWidth = 4.17

# --- View: zoomed to 160%, with E10 selected ---
$excel.ActiveWindow.Zoom = 160
$ws.Range("E10").Select()
